# PaintingZigeunerpaar.xlsx fix:
#  - rename label "Objekt-ID" -> "ObjektID" (cell A2)
#  - give the data rows (B2 and A3:B32) a new monospace font
#  - move the active selection to A2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the property label text in A2 ("Objekt-ID" -> "ObjektID")
$ws.Range("A2").Value = "ObjektID"

# 2. Give B2 the new monospace font: Liberation Mono / Courier New /
#    DejaVu Sans Mono / Lucida Sans Typewriter, size 10, family 3 (modern/monospace).
$src = $ws.Range("B2")
$src.Font.Name = "Liberation Mono;Courier New;DejaVu Sans Mono;Lucida Sans Typewriter"
$src.Font.Size = 10
$src.Font.Family = 3

# Propagate the same formatting to the rest of the data rows (A3:B32) by copying
# B2's now-updated format over them, so they all share the same cell style
# instead of each creating its own.
$src.Copy()
$dst = $ws.Range("A3:B32")
$dst.PasteSpecial(-4122) | Out-Null

# 3. Move the active cell selection from B32 to A2
$ws.Range("A2").Select() | Out-Null
